$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Tasks 03-18 to 03-25")
$dst = $wb.Worksheets.Item("Tasks 03-11 to 03-18")
$src.Range("F2").Copy()
$dst.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats
Write-Host "dest format copied:" $dst.Range("Z1").Interior.Color
